$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two existing hyperlinks up front (Range.Hyperlinks.Delete() is the
# only reliable removal path in this host - Hyperlinks.Item(n).Delete() is a no-op).
$ws.Range("C3").Hyperlinks.Delete()

# Insert a new column before column C ("Resource"), shifting the rest right.
$ws.Columns("C").Insert()

# Populate the new "Resource" column.
$ws.Range("C1").Value = "Resource"
$ws.Range("C2").Value = "pet"
$ws.Range("C3").Value = "pet"

# Re-point the conditional formatting rule that used to live on G2 to its new
# location H2 (both the applies-to range and the formula text reference it).
$fc = $ws.Range("G2").FormatConditions.Item(1)
$fc.Formula1 = "=LEN(TRIM(H2))>0"
$fc.ModifyAppliesToRange($ws.Range("H2"))

# Re-create the two hyperlinks at their new (shifted) locations: D3 and D4.
# D4 has no backing cell value in the original sheet (it was an orphan
# hyperlink entry), so add it, then delete the helper row it creates.
$ws.Hyperlinks.Add($ws.Range("D3"), "https://live.virtualandemo.com/api/pets/findByTags?tags=grey")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://live.virtualandemo.com/api/pets/findByTags?tags=grey", $null, $null, "https://live.virtualandemo.com/api/pets/findByTags?tags=grey")
$ws.Rows(4).Delete()

# Hyperlinks.Add() stamps the built-in "Hyperlink" style; restore D3's
# original font formatting (bold underlined Inconsolata in #1155CC).
$f3 = $ws.Range("D3").Font
$f3.Name = "Inconsolata"
$f3.Bold = $true
$f3.Underline = 2
$f3.Size = 11
$f3.Color = 13391121

# Match the post-edit selection (C4).
$ws.Range("C4").Select()

Write-Output "done"
